# Update "想去人数" (want-to-go count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 302
$ws1.Range("F4").Value = 1210
$ws1.Range("F5").Value = 608

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 302
$ws4.Range("F4").Value = 1210
$ws4.Range("F6").Value = 608
